$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 793817.9
$ws.Range("I2").Value = 1010278.25
$ws.Range("J2").Value = 129.66667
$ws.Range("K2").Value = 1010278.25
$ws.Range("L2").Value = 129.66667
$ws.Range("M2").Value = -1010165.25
$ws.Range("N2").Value = -355.66667
$ws.Range("H17").Value = 971.5238000000001
$ws.Range("J17").Value = 971.5238000000001
$ws.Range("L17").Value = 2914.5714
$ws.Range("N17").Value = -3250.5714
$ws.Range("H69").Value = 4433.3335
$ws.Range("J69").Value = 4433.3335
$ws.Range("L69").Value = 13300.0005
$ws.Range("N69").Value = -15048.0005
$ws.Range("H72").Value = 4433.3335
$ws.Range("J72").Value = 4433.3335
$ws.Range("L72").Value = 39900.0015
$ws.Range("N72").Value = -48636.0015
$ws.Range("H106").Value = 3166.6667
$ws.Range("I106").Value = 3166.6667
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3166.6667
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2535.6667
$ws.Range("N106").ClearContents()
$ws.Range("H124").Value = 51000
$ws.Range("J124").Value = 51000
$ws.Range("L124").Value = 51000
$ws.Range("N124").Value = -60820

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7281.323
$ws.Range("I32").Value = 4360.829
$ws.Range("J32").Value = 16931.652
$ws.Range("K32").Value = 4360.829
$ws.Range("L32").Value = 16931.652
$ws.Range("M32").Value = -4073.829
$ws.Range("N32").Value = -17505.652

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 22030
$ws.Range("J31").Value = 22030
$ws.Range("L31").Value = 22030
$ws.Range("N31").Value = -22534
$ws.Range("H86").Value = 1462.826
$ws.Range("I86").Value = 1374.2142
$ws.Range("J86").Value = 1600.6666
$ws.Range("K86").Value = 1374.2142
$ws.Range("L86").Value = 1600.6666
$ws.Range("M86").Value = -251.2141999999999
$ws.Range("N86").Value = -3846.6666
$ws.Range("H89").Value = 1462.826
$ws.Range("I89").Value = 1374.2142
$ws.Range("J89").Value = 1600.6666
$ws.Range("K89").Value = 6871.071
$ws.Range("L89").Value = 8003.333000000001
$ws.Range("M89").Value = -1255.071
$ws.Range("N89").Value = -19235.333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2348.6487
$ws.Range("I31").Value = 1499.5333
$ws.Range("J31").Value = 5987.7144
$ws.Range("K31").Value = 1499.5333
$ws.Range("L31").Value = 5987.7144
$ws.Range("M31").Value = -1204.5333
$ws.Range("N31").Value = -6577.7144
$ws.Range("H34").Value = 2348.6487
$ws.Range("I34").Value = 1499.5333
$ws.Range("J34").Value = 5987.7144
$ws.Range("K34").Value = 1499.5333
$ws.Range("L34").Value = 5987.7144
$ws.Range("M34").Value = -1297.5333
$ws.Range("N34").Value = -6391.7144
$ws.Range("H99").Value = 1953.9722
$ws.Range("I99").Value = 1867.4231
$ws.Range("J99").Value = 2179
$ws.Range("K99").Value = 1867.4231
$ws.Range("L99").Value = 2179
$ws.Range("M99").Value = -369.4231
$ws.Range("N99").Value = -5175
$ws.Range("H126").Value = 1953.9722
$ws.Range("I126").Value = 1867.4231
$ws.Range("J126").Value = 2179
$ws.Range("K126").Value = 5602.2693
$ws.Range("L126").Value = 6537
$ws.Range("M126").Value = -3132.2693
$ws.Range("N126").Value = -11477

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 774.5282999999999
$ws.Range("I131").Value = 348.86667
$ws.Range("J131").Value = 942.5526
$ws.Range("K131").Value = 1046.60001
$ws.Range("L131").Value = 2827.6578
$ws.Range("M131").Value = 3993.39999
$ws.Range("N131").Value = -12907.6578
$ws.Range("H134").Value = 167735.89
$ws.Range("I134").Value = 527678.4
$ws.Range("J134").Value = 4904.7617
$ws.Range("K134").Value = 1583035.2
$ws.Range("L134").Value = 14714.2851
$ws.Range("M134").Value = -1577965.2
$ws.Range("N134").Value = -24854.2851

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2848.182
$ws.Range("I132").Value = 2698.9465
$ws.Range("K132").Value = 8096.8395
$ws.Range("M132").Value = -5566.8395

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 440.3684
$ws.Range("I22").Value = 315.41666
$ws.Range("J22").Value = 654.5714
$ws.Range("K22").Value = 315.41666
$ws.Range("L22").Value = 654.5714
$ws.Range("M22").Value = -20.41665999999998
$ws.Range("N22").Value = -1244.5714
$ws.Range("H27").Value = 440.3684
$ws.Range("I27").Value = 315.41666
$ws.Range("J27").Value = 654.5714
$ws.Range("K27").Value = 315.41666
$ws.Range("L27").Value = 654.5714
$ws.Range("M27").Value = -208.41666
$ws.Range("N27").Value = -868.5714
$ws.Range("H46").Value = 1111
$ws.Range("I46").Value = 999.7143
$ws.Range("J46").Value = 1890
$ws.Range("K46").Value = 999.7143
$ws.Range("L46").Value = 1890
$ws.Range("M46").Value = -811.7143
$ws.Range("N46").Value = -2266
$ws.Range("H55").Value = 215.29033
$ws.Range("I55").Value = 162.76923
$ws.Range("J55").Value = 488.4
$ws.Range("K55").Value = 162.76923
$ws.Range("L55").Value = 488.4
$ws.Range("M55").Value = 10.23077000000001
$ws.Range("N55").Value = -834.4
$ws.Range("H87").Value = 13189
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 13189
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 13189
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -15435
$ws.Range("H90").Value = 13189
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 13189
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 39567
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -50799
$ws.Range("H122").Value = 5368.282
$ws.Range("I122").Value = 5671.4546
$ws.Range("K122").Value = 17014.3638
$ws.Range("M122").Value = -14564.3638
$ws.Range("H132").Value = 1234.3906
$ws.Range("I132").Value = 905.6226
$ws.Range("J132").Value = 2818.4546
$ws.Range("K132").Value = 2716.8678
$ws.Range("L132").Value = 8455.363799999999
$ws.Range("M132").Value = -186.8678
$ws.Range("N132").Value = -13515.3638

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 41260
$ws.Range("I95").Value = 41000
$ws.Range("J95").Value = 41325
$ws.Range("K95").Value = 41000
$ws.Range("L95").Value = 41325
$ws.Range("M95").Value = -38254
$ws.Range("N95").Value = -46817
$ws.Range("H132").Value = 1059.409
$ws.Range("I132").Value = 516.70966
$ws.Range("J132").Value = 2353.5386
$ws.Range("K132").Value = 1550.12898
$ws.Range("L132").Value = 7060.6158
$ws.Range("M132").Value = 979.87102
$ws.Range("N132").Value = -12120.6158
$ws.Range("H136").Value = 4560.3335
$ws.Range("I136").Value = 1188.25
$ws.Range("J136").Value = 11304.5
$ws.Range("K136").Value = 3564.75
$ws.Range("L136").Value = 33913.5
$ws.Range("M136").Value = -1014.75
$ws.Range("N136").Value = -39013.5
